$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose column-A "match code" label (e.g. "M001") is a duplicate of the
# code already shown on the first row of its match-block. The transaction
# block logic should only stamp the code on the first row of each block, so
# clear column A on every subsequent row of each block.
$rowsToClear = @(
    54, 55, 57, 58, 59, 60, 61, 62, 69, 70, 72, 73, 74, 87, 88,
    90, 91, 92, 108, 109, 110, 230, 231, 232, 236, 237, 238, 256, 257, 258,
    277, 278, 279, 338, 339, 341, 342, 344, 345, 347, 348, 349, 362, 363, 365,
    366, 368, 369, 370, 380, 381, 382, 383, 384, 433, 434, 435, 442, 443, 444,
    462, 463, 464, 468, 469, 471, 472, 474, 475, 476, 477, 478, 509, 510, 511,
    521, 522, 523, 571, 572, 574, 575, 576, 584, 585, 586, 587, 589, 590, 592,
    593, 594, 641, 642, 644, 645, 647, 648, 649, 671, 672, 674, 675, 676, 677,
    678, 682, 683, 684, 685, 686, 696, 697, 698, 702, 703, 704, 711, 712, 714,
    715, 717, 718, 719, 723, 724, 725, 737, 738, 740, 741, 743, 744, 745, 763,
    764, 766, 767, 768, 772, 773, 774, 921, 922, 924, 925, 927, 928, 929, 942,
    943, 944, 960, 961, 962, 963, 964, 973, 974, 975, 979, 980, 981, 982, 983,
    985, 986, 987, 1011, 1012, 1013, 1032, 1033, 1034, 1101, 1102, 1104, 1105, 1106, 1127,
    1128, 1129, 1130, 1131, 1132, 1134, 1135, 1136, 1161, 1162, 1163, 1170, 1171, 1172, 1197,
    1198, 1199, 1220, 1221, 1222
)

foreach ($r in $rowsToClear) {
    $ws.Cells.Item($r, 1).Value = ""
}
